# Updated to pull LA and ICB data from FT for HC indicator (ID = 71)
#
# - FT_indicators!C19 (IndicatorID 71, FingerTips_id 91041) changes its
#   AreaType from "England" to "LA". This makes the previously-referenced
#   "England" shared string unused, so it drops out of sharedStrings.xml.
# - The workbook's active sheet moves from "meta_only" to "FT_indicators",
#   and the cursor/selection on each sheet is updated to reflect where the
#   author was last working (meta_only: H28, FT_indicators: E18).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("meta_only")
$wsFT   = $wb.Worksheets.Item("FT_indicators")

# Capture the final cursor position left on the sheet that is no longer
# active before we switch focus away from it.
$wsMeta.Range("H28").Select()

# Apply the actual data change: AreaType for IndicatorID 71 goes from
# "England" to "LA".
$wsFT.Range("C19").Value = "LA"

# Leave FT_indicators as the active sheet/tab with E18 selected.
$wsFT.Activate()
$wsFT.Range("E18").Select()
